$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.827.24"
$ws.Range("E2").Value = "  +0.93%  "

# Row 3
$ws.Range("D3").Value = "3.792.16"
$ws.Range("E3").Value = "  +0.45%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "444.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.54%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.75%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.55%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.735"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.55%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.76%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000318"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.12%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.96%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.18%  "

# Row 14
$ws.Range("D14").Value = "4.446.94"
$ws.Range("E14").Value = "  +0.34%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -9.22%  "

# Row 16
$ws.Range("D16").Value = "3.825.61"
$ws.Range("E16").Value = "  +0.32%  "

# Row 17
$ws.Range("E17").Value = "  -0.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.13%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.91%  "

# Row 20
$ws.Range("D20").Value = "66.989.23"
$ws.Range("E20").Value = "  +0.69%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "417.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.84%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.31%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.37%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.51%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "37.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.12%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.37%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.30%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +21.92%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "730.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.59%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.08%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.132"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.22%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.96%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +16.06%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.157"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.22%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.65%  "

# Row 37
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.06%  "

# Row 38
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +21.26%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0473"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.07%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.08%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.16%  "

# Row 42
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0671"
$ws.Range("E42").Value = "  -10.94%  "

# Row 43
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.139"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.55%  "

# Row 44
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.332"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +15.21%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "27.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.80%  "

# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.95%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +18.06%  "

# Row 48
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.49%  "

# Row 49
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.25%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.21%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.81%  "
